# Applies the automatic-update diff to the "Artfynd" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple "Taxonsorteringsordning" (column B) bumps ---------------------
$ws.Range("B6").Value = 83216
$ws.Range("B7").Value = 80384

$ws.Range("B11").Value = 91772
$ws.Range("B12").Value = 83207
$ws.Range("B13").Value = 83207
$ws.Range("B14").Value = 83216
$ws.Range("B15").Value = 91772
$ws.Range("B16").Value = 91823
$ws.Range("B17").Value = 91809

$ws.Range("B24").Value = 83090
$ws.Range("B25").Value = 91772
$ws.Range("B26").Value = 91772
$ws.Range("B27").Value = 83216
$ws.Range("B28").Value = 83216
$ws.Range("B29").Value = 92228
$ws.Range("B30").Value = 80350
$ws.Range("B32").Value = 91809
$ws.Range("B33").Value = 83216
$ws.Range("B34").Value = 83090

# --- Row 9 / Row 10 content swap (with K/L/M/N/AC moving to row 9) --------
$ws.Range("A9").Value = 131066770
$ws.Range("B9").Value = 57884
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = ""
$ws.Range("M9").Value = "äldre spår"
$ws.Range("N9").Value = ""
$ws.Range("Q9").Value = 425323
$ws.Range("R9").Value = 6712206
$ws.Range("AC9").Value = "Ringhack på gran"

$ws.Range("A10").Value = 131066783
$ws.Range("B10").Value = 83090
$ws.Range("E10").Value = 1312
$ws.Range("F10").Value = "Gammelgransskål"
$ws.Range("G10").Value = "Pseudographis pinicola"
$ws.Range("H10").Value = "(Nyl.) Rehm"
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = ""
$ws.Range("Q10").Value = 425170
$ws.Range("R10").Value = 6712292
$ws.Range("AC10").Value = ""

# --- Row 21 / Row 22 content swap ------------------------------------------
$ws.Range("A21").Value = 131066778
$ws.Range("B21").Value = 81229
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 1049
$ws.Range("F21").Value = "Kortskaftad ärgspik"
$ws.Range("G21").Value = "Microcalicium ahlneri"
$ws.Range("H21").Value = "Tibell"
$ws.Range("Q21").Value = 425336
$ws.Range("R21").Value = 6712202

$ws.Range("A22").Value = 131066766
$ws.Range("B22").Value = 92180
$ws.Range("D22").Value = "VU"
$ws.Range("E22").Value = 2062
$ws.Range("F22").Value = "Ulltickeporing"
$ws.Range("G22").Value = "Skeletocutis brevispora"
$ws.Range("H22").Value = "Niemelä"
$ws.Range("Q22").Value = 425069
$ws.Range("R22").Value = 6712285
